$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are free-form text (e.g. "30.630.28"); force Text format
# so Excel does not silently reinterpret/round them as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.625.96"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.138.81"
$ws.Range("E3").Value = "  +1.46%  "
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.62"
$ws.Range("E5").Value = "  +5.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.009"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5269"
$ws.Range("E7").Value = "  +0.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4554"
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.78"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09133"
$ws.Range("E10").Value = "  +2.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.183"
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.81"
$ws.Range("E12").Value = "  +2.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.142.18"
$ws.Range("E13").Value = "  +1.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.867"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.131"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("E16").Value = "  +5.58%  "
$ws.Range("E17").Value = "  +2.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.011"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06711"
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.49"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.009"
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.345"
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.695.91"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.81"
$ws.Range("E24").Value = "  +3.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.383"
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.372.44"
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.56"
$ws.Range("E27").Value = "  +1.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.649"
$ws.Range("E28").Value = "  +4.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "164.80"
$ws.Range("E29").Value = "  +1.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "137.03"
$ws.Range("E30").Value = "  +2.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.221"
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("E33").Value = "  +1.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.359"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.015"
$ws.Range("E35").Value = "  +1.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.196"
$ws.Range("E36").Value = "  +7.46%  "
$ws.Range("E37").Value = "  +1.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02654"
$ws.Range("E38").Value = "  +2.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06923"
$ws.Range("E39").Value = "  +1.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2332"
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6979"
$ws.Range("E42").Value = "  +1.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.273"
$ws.Range("E43").Value = "  +2.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.79"
$ws.Range("E44").Value = "  +5.08%  "
$ws.Range("E45").Value = "  +1.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6459"
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000370"
$ws.Range("E47").Value = "  +5.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.756"
$ws.Range("E48").Value = "  +2.53%  "
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.02"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07297"
$ws.Range("E51").Value = "  +2.27%  "

Write-Output "Updated cryptos list"